# Friday commit: pull config references out of the job report to ease
# distribution of the workflow without the RE framework being required.
#
# The "ReportResults" sheet used to list six hard-coded test-run rows
# (cols A + D) starting at row 2. Those six rows get pushed down to
# rows 24-29 (D28's transaction number also gets corrected), and the
# freed-up rows 2-13 are filled with the validation-message test cases
# that used to live in the RE framework's config.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReportResults")

# --- Capture the six existing test rows (cols A + D) before moving them ---
$oldRows = @()
for ($r = 2; $r -le 7; $r++) {
    $oldRows += ,@($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 4).Value2)
}

# --- Clear rows 2:7 so the sheet can be rebuilt cleanly ---
$ws.Range("A2:G7").Clear()

# --- New rows 2-13: validation-message test cases pulled out of config.
#     Written in this order so the workbook's shared-string table builds
#     up the same way it did for the author (first-seen order below). ---
$ws.Cells.Item(2, 1).Value2 = "No RM number found in email"
$ws.Cells.Item(3, 1).Value2 = "RM number not a 10 digit code"
$ws.Cells.Item(7, 1).Value2 = "Extra digits: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(10, 1).Value2 = "Bad char: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(5, 1).Value2 = "Blank Field: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(13, 1).Value2 = "Valid form with 60 character short description"
$ws.Cells.Item(4, 1).Value2 = "Blank Field: in OUC: Fac, Sec. UA: Opal, Fund"

$ws.Cells.Item(6, 1).Value2 = "Blank Field: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(8, 1).Value2 = "Extra digits: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(9, 1).Value2 = "Extra digits: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(11, 1).Value2 = "Bad char: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"
$ws.Cells.Item(12, 1).Value2 = "Bad char: in OUC: Fac, Sch, Sec, Site. UA: Opal, Site, Fund, Func"

# --- Re-write the original six rows further down (rows 24-29). Row 28's
#     (old row 6) transaction number is corrected: 2019003060 -> 2018001996 ---
$oldRows[4][1] = 2018001996

$row = 24
foreach ($pair in $oldRows) {
    $ws.Cells.Item($row, 1).Value2 = $pair[0]
    $ws.Cells.Item($row, 4).Value2 = $pair[1]
    $row++
}

# --- Column widths Excel recorded on save (now narrower, auto-fit-style) ---
$ws.Columns.Item(1).ColumnWidth = 20.42578125
$ws.Columns.Item(2).ColumnWidth = 13.7109375
$ws.Columns.Item(3).ColumnWidth = 10.5703125
$ws.Columns.Item(4).ColumnWidth = 21
$ws.Columns.Item(5).ColumnWidth = 12.140625
$ws.Columns.Item(6).ColumnWidth = 13.7109375

# --- Selection moved to A14 when the author saved, but SummaryResults
#     stayed the active tab, so re-activate it afterwards ---
$ws.Range("A14").Select()
$wb.Worksheets.Item("SummaryResults").Activate()
